# Dokumentáció.pptx - "Add files via upload" edit
#
# On slide 3 (Csapattagok feladata), shape "Tartalom helye 2" has a
# paragraph reading:
#   "                                - Képek, információk gyűjtése  "
# It gets split into two runs: a brand-new leading run that holds
# nothing but spaces, and the original run with its leading whitespace
# stripped off so it now starts directly with the dash.

$needle = "Képek, információk gyűjtése"

$p = $ppt.ActivePresentation

foreach ($s in $p.Slides) {
    foreach ($sh in $s.Shapes) {
        if (-not $sh.HasTextFrame) { continue }

        $tr = $sh.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count

        for ($i = 1; $i -le $paraCount; $i++) {
            $para = $tr.Paragraphs($i, 1)

            if ($para.Text.Contains($needle)) {
                $run = $para.Runs(1, 1)

                # Insert a brand-new run made up solely of spaces right
                # before the existing run (41 spaces total).
                $run.InsertBefore("                                         ")

                # The original run is now the second run of the
                # paragraph; strip its leading spaces so it starts
                # directly with the dash.
                $run2 = $para.Runs(2, 1)
                $run2.Text = "- Képek, információk gyűjtése  "
            }
        }
    }
}
